# Update with Correct Forecast output
#
# - Rename Sheet1 to "Sales vs PO".
# - Insert an "Order Week" column (the original order date) ahead of the
#   PO_Requested_Qty column; the sale/forecast date column ("ds") is shifted
#   forward one week and PO_Requested_Qty on this sheet is zeroed out (it now
#   lives on its own sheet).
# - Add "Weekly Growth": the real weekly PO quantities plus week-over-week
#   growth %.
# - Add "Volume Insights": aggregate PO stats (total/average/max/min).
# - Add "Prediction Info": the forecast for next week's PO quantity (simple
#   linear trend extrapolation over the weekly quantities).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Original data captured from Sheet1 before the edit: (ds serial, y, PO qty)
# ---------------------------------------------------------------------------
$origDates = @(45327,45334,45341,45348,45355,45362,45369,45376,45383,45390,45397,45404,45411,45418,45425,45432,45439,45446,45453,45460,45467,45474,45481,45488,45495,45502,45509,45516,45523,45530,45537,45544,45551,45558,45565,45572,45579,45586,45593,45600,45607,45614,45621,45628,45635,45642,45649)
$origY     = @(0,0,1,3,5,4,4,4,6,3,3,3,7,7,3,7,5,23,2,0,0,2,3,0,2,1,3,3,83,1,1,1,1,0,0,3,0,1,1,0,5,2,3,2,1,3,0)
$origPO    = @(0,72,48,24,16,0,8,0,32,0,16,0,0,0,32,0,16,16,112,0,16,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$n = $origDates.Length

# Non-zero weekly PO quantities (and their dates) -- these move to the new
# "Weekly Growth" sheet.
$growthDates = @()
$growthQty   = @()
for ($i = 0; $i -lt $n; $i++) {
    if ($origPO[$i] -ne 0) {
        $growthDates += $origDates[$i]
        $growthQty   += $origPO[$i]
    }
}
$gn = $growthDates.Length

# ---------------------------------------------------------------------------
# Clone Sheet1 three times (before mutating it) so the new sheets inherit the
# same sheetPr / page setup / base styles (header style, date style, ...).
# ---------------------------------------------------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws1.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)

# ===========================================================================
# Sheet1 -> "Sales vs PO"
# ===========================================================================
$ws1.Name = "Sales vs PO"

# Insert a new column ahead of the existing PO_Requested_Qty column (col C)
$ws1.Columns.Item(3).Insert()

# Header row (the newly-inserted column C cell already inherited the bold
# header style s=1 from its neighbours via the column insert)
$ws1.Cells.Item(1,1).Value2 = "ds"
$ws1.Cells.Item(1,2).Value2 = "y"
$ws1.Cells.Item(1,3).Value2 = "Order Week"
$ws1.Cells.Item(1,4).Value2 = "PO_Requested_Qty"

for ($i = 0; $i -lt $n; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r,1).Value2 = $origDates[$i] + 6
    $ws1.Cells.Item($r,2).Value2 = $origY[$i]
    $ws1.Cells.Item($r,3).Value2 = $origDates[$i]
    $ws1.Cells.Item($r,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws1.Cells.Item($r,4).Value2 = 0
}

# ===========================================================================
# Sheet2 -> "Weekly Growth"
# ===========================================================================
$ws2.Name = "Weekly Growth"
$ws2.Cells.Clear()

$ws2.Cells.Item(1,1).Value2 = "ds"
$ws2.Cells.Item(1,2).Value2 = "PO_Requested_Qty"
$ws2.Cells.Item(1,3).Value2 = "Growth%"

for ($i = 0; $i -lt $gn; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r,1).Value2 = $growthDates[$i]
    $ws2.Cells.Item($r,2).Value2 = $growthQty[$i]
    if ($i -eq 0) {
        $ws2.Cells.Item($r,3).Value2 = 0
    } else {
        $prev = $growthQty[$i-1]
        $cur  = $growthQty[$i]
        $ws2.Cells.Item($r,3).Value2 = ($cur / $prev * 100) - 100
    }
}

# Re-apply the header style and the date-column style (grabbed from the
# still-pristine-looking "Sales vs PO" header / date cells) rather than the
# default styles left behind by Cells.Clear().
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

$lastRow2 = $gn + 1
$ws1.Range("A2:A$lastRow2").Copy()
$ws2.Range("A2:A$lastRow2").PasteSpecial(-4122)

# ===========================================================================
# Sheet3 -> "Volume Insights"
# ===========================================================================
$ws3.Name = "Volume Insights"
$ws3.Cells.Clear()

$total = 0
$max = $growthQty[0]
$min = $growthQty[0]
foreach ($q in $growthQty) {
    $total += $q
    if ($q -gt $max) { $max = $q }
    if ($q -lt $min) { $min = $q }
}
$avg = $total / $gn

$ws3.Cells.Item(1,1).Value2 = "Total_PO_Quantity"
$ws3.Cells.Item(1,2).Value2 = "Average_PO_Quantity"
$ws3.Cells.Item(1,3).Value2 = "Max_PO_Quantity"
$ws3.Cells.Item(1,4).Value2 = "Min_PO_Quantity"

$ws3.Cells.Item(2,1).Value2 = $total
$ws3.Cells.Item(2,2).Value2 = $avg
$ws3.Cells.Item(2,3).Value2 = $max
$ws3.Cells.Item(2,4).Value2 = $min

$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# ===========================================================================
# Sheet4 -> "Prediction Info"
# ===========================================================================
$ws4.Name = "Prediction Info"
$ws4.Cells.Clear()

# Simple linear regression (least squares) over the weekly PO quantities,
# extrapolated one point beyond the observed series -> next week's forecast.
$sumX = 0
$sumY = 0
$sumXY = 0
$sumXX = 0
for ($i = 0; $i -lt $gn; $i++) {
    $x = $i
    $y = $growthQty[$i]
    $sumX += $x
    $sumY += $y
    $sumXY += $x * $y
    $sumXX += $x * $x
}
$slope = (($gn * $sumXY) - ($sumX * $sumY)) / (($gn * $sumXX) - ($sumX * $sumX))
$intercept = ($sumY - ($slope * $sumX)) / $gn
$prediction = ($slope * $gn) + $intercept

$ws4.Cells.Item(1,1).Value2 = "Predicted_Next_Week_PO_Quantity"
$ws4.Cells.Item(2,1).Value2 = $prediction

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Restore the first sheet as the active / selected sheet
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
